$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) to remain plain text so numeric-looking strings
# (e.g. "1.002") are not silently converted to real numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "26.938.65"
$ws.Range("E2").Value = "  +0.05%  "

# Row 3
$ws.Range("D3").Value = "1.817.82"
$ws.Range("E3").Value = "  +0.43%  "

# Row 4
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.21%  "

# Row 5
$ws.Range("D5").Value = "309.68"
$ws.Range("E5").Value = "  -0.16%  "

# Row 6
$ws.Range("E6").Value = "  +0.19%  "

# Row 7
$ws.Range("D7").Value = "0.4655"
$ws.Range("E7").Value = "  +0.67%  "

# Row 8
$ws.Range("D8").Value = "0.3665"
$ws.Range("E8").Value = "  -1.29%  "

# Row 9
$ws.Range("D9").Value = "0.07368"
$ws.Range("E9").Value = "  -0.11%  "

# Row 10
$ws.Range("D10").Value = "0.8724"
$ws.Range("E10").Value = "  -0.53%  "

# Row 11
$ws.Range("D11").Value = "20.27"
$ws.Range("E11").Value = "  -1.16%  "

# Row 12
$ws.Range("D12").Value = "1.815.95"
$ws.Range("E12").Value = "  -3.56%  "

# Row 13
$ws.Range("E13").Value = "  +0.53%  "

# Row 14
$ws.Range("D14").Value = "0.07109"
$ws.Range("E14").Value = "  +0.97%  "

# Row 15
$ws.Range("D15").Value = "6.509"
$ws.Range("E15").Value = "  -0.10%  "

# Row 16
$ws.Range("D16").Value = "91.32"
$ws.Range("E16").Value = "  -1.27%  "

# Row 17
$ws.Range("D17").Value = "1.004"
$ws.Range("E17").Value = "  +0.33%  "

# Row 18
$ws.Range("E18").Value = "  -0.33%  "

# Row 20
$ws.Range("E20").Value = "  -0.59%  "

# Row 21
$ws.Range("D21").Value = "26.963.93"
$ws.Range("E21").Value = "  +0.08%  "

# Row 22
$ws.Range("E22").Value = "  -0.63%  "

# Row 23
$ws.Range("D23").Value = "10.59"
$ws.Range("E23").Value = "  -0.43%  "

# Row 24
$ws.Range("D24").Value = "2.036.69"
$ws.Range("E24").Value = "  -1.86%  "

# Row 25
$ws.Range("D25").Value = "1.895"
$ws.Range("E25").Value = "  -0.24%  "

# Row 26
$ws.Range("D26").Value = "151.24"
$ws.Range("E26").Value = "  -0.11%  "

# Row 27
$ws.Range("D27").Value = "18.45"
$ws.Range("E27").Value = "  +0.30%  "

# Row 28
$ws.Range("D28").Value = "2.136"
$ws.Range("E28").Value = "  -1.12%  "

# Row 29
$ws.Range("D29").Value = "5.253"
$ws.Range("E29").Value = "  -1.80%  "

# Row 30
$ws.Range("D30").Value = "116.80"
$ws.Range("E30").Value = "  +0.72%  "

# Row 31
$ws.Range("D31").Value = "0.08897"
$ws.Range("E31").Value = "  -0.04%  "

# Row 32
$ws.Range("D32").Value = "0.7588"
$ws.Range("E32").Value = "  +0.38%  "

# Row 33
$ws.Range("D33").Value = "1.164"
$ws.Range("E33").Value = "  +0.35%  "

# Row 34
$ws.Range("D34").Value = "4.484"
$ws.Range("E34").Value = "  +0.64%  "

# Row 35
$ws.Range("E35").Value = "  -0.35%  "

# Row 36
$ws.Range("E36").Value = "  +0.20%  "

# Row 37
$ws.Range("D37").Value = "1.095"
$ws.Range("E37").Value = "  -0.81%  "

# Row 38
$ws.Range("D38").Value = "0.05289"
$ws.Range("E38").Value = "  +0.68%  "

# Row 39
$ws.Range("D39").Value = "0.01946"
$ws.Range("E39").Value = "  -1.32%  "

# Row 40
$ws.Range("D40").Value = "2.981"
$ws.Range("E40").Value = "  +1.88%  "

# Row 41
$ws.Range("D41").Value = "0.5292"
$ws.Range("E41").Value = "  -0.81%  "

# Row 42
$ws.Range("D42").Value = "7.147"
$ws.Range("E42").Value = "  -1.15%  "

# Row 43
$ws.Range("D43").Value = "2.325"
$ws.Range("E43").Value = "  -3.85%  "

# Row 45
$ws.Range("D45").Value = "8.425"
$ws.Range("E45").Value = "  -1.10%  "

# Row 46
$ws.Range("D46").Value = "0.4851"
$ws.Range("E46").Value = "  -2.76%  "

# Row 47
$ws.Range("D47").Value = "10.41"
$ws.Range("E47").Value = "  +0.61%  "

# Row 48
$ws.Range("E48").Value = "  +0.25%  "

# Row 49
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "1.661"
$ws.Range("E49").Value = "  -0.82%  "

# Row 50
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "103.23"
$ws.Range("E50").Value = "  -0.70%  "

# Row 51
$ws.Range("D51").Value = "0.06289"
$ws.Range("E51").Value = "  -0.05%  "
